$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on D/E columns so numeric-looking strings
# (e.g. "10.60", "0.0361") are preserved exactly as text, matching the
# source workbook convention of storing Price/Volume columns as inline strings.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.119.51'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.08%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.053.18'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.12%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.28'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.55%  '

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.85%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.79%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.053.24'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.19%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.74%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.82'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.29%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.447'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.76%  '

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.40%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.21'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.19%  '

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.83%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.556.54'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.21%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.13'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.82%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.088.29'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.08%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.051.09'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.21%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '478.07'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.01%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.33%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.44%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.24%  '

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.62%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.43'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.63%  '

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.98%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.60'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +6.35%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.31%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.35'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.15%  '

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.24%  '

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.07%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.03%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.58'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.51%  '

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.52%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.39%  '

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.52%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.54%  '

$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.89'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.32%  '

$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.21'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.27%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.23'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.41%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.27%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '433.92'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.99%  '

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.18%  '

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.03%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0361'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.828.80'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.28%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '38.35'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.48%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.05'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.56%  '

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.02%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.18'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.38%  '

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.58%  '
